# This script updates the "quiz" marksheet worksheet so that it reflects an
# actual graded attempt instead of the original "nothing attempted / Absent"
# placeholder state, and collapses the sheet from three parallel
# student-answer blocks (A/B, D/E, G/H) down to a single, correctly scored
# block (mirrors the fix for handling float/locale input described in the
# commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Summary table (rows 10-12)
# ---------------------------------------------------------------------
# Give the "No./Marking/Total" row labels in column A the same header
# style used elsewhere (mtitleStyle), matching the rest of the sheet.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# Row 10 ("No.") - counts of right / wrong / not-attempted / max questions
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 28

# Row 11 ("Marking") - marks awarded per right/wrong answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 ("Total") - right*mark, wrong*mark, and the "score/max" summary
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "56/112"

# ---------------------------------------------------------------------
# 2. Remove the third ("G/H") answer block entirely
# ---------------------------------------------------------------------
$ws.Columns("G:H").Delete()

# ---------------------------------------------------------------------
# 3. Trim the second ("D/E") answer block down to just the first three
#    questions (rows 16-18); the rest is no longer needed.
# ---------------------------------------------------------------------
$ws.Range("D19:E40").ClearContents()

# Give the remaining D/E student-answer cells (rows 16-18) their graded
# (correct) styling, matching the correct-answer values already in E16:E18.
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------
# 4. Fill in the "Student Ans" column (A16:A40) for the main answer block
#    with the actual answer given, colored green (correctStyle) when it
#    matches the "Correct Ans" column B, red (incorrectStyle) when it
#    doesn't, and left blank with normalStyle when the question was not
#    attempted.
# ---------------------------------------------------------------------
$answers = @{
    16 = @("correctStyle",   "Option A")
    17 = @("correctStyle",   "Option D")
    18 = @("correctStyle",   "Option B")
    19 = @("incorrectStyle", "Option D")
    20 = @("correctStyle",   "Option B")
    21 = @("incorrectStyle", "Option B")
    22 = @("incorrectStyle", "Option A")
    23 = @("correctStyle",   "Option D")
    24 = @("normalStyle",    $null)
    25 = @("incorrectStyle", "Option D")
    26 = @("correctStyle",   "Option C")
    27 = @("normalStyle",    $null)
    28 = @("incorrectStyle", "Option B")
    29 = @("correctStyle",   "Option D")
    30 = @("correctStyle",   "Option B")
    31 = @("incorrectStyle", "Option C")
    32 = @("correctStyle",   "Option C")
    33 = @("correctStyle",   "Option D")
    34 = @("normalStyle",    $null)
    35 = @("incorrectStyle", "Option B")
    36 = @("correctStyle",   "Option A")
    37 = @("normalStyle",    $null)
    38 = @("correctStyle",   "Option A")
    39 = @("correctStyle",   "Option D")
    40 = @("incorrectStyle", "Option B")
}

foreach ($row in $answers.Keys) {
    $info = $answers[$row]
    $style = $info[0]
    $value = $info[1]
    $cell = $ws.Range("A$row")
    $cell.Style = $style
    if ($null -ne $value) {
        $cell.Value = $value
    }
}
